$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 49272

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2634

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2288

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1219
